# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# D-column prices that are plain numeric-looking strings (e.g. "211.48") are written
# with a leading apostrophe so Excel keeps them as literal text instead of silently
# converting them to floating-point numbers (this mirrors how the source data keeps
# thousand-dot-separated prices like "26.714.11" as text already).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '26.714.11'
$ws.Range("E2").Value = '  +0.35%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.599.50'
$ws.Range("E3").Value = '  +0.32%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.20%  '

# Row 5: BNB
$ws.Range("D5").Value = '''211.48'
$ws.Range("E5").Value = '  +0.14%  '

# Row 6: XRP
$ws.Range("D6").Value = '''0.513'
$ws.Range("E6").Value = '  -0.52%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.19%  '

# Row 8: Dogecoin
$ws.Range("E8").Value = '  +0.44%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  +1.18%  '

# Row 10: Solana
$ws.Range("D10").Value = '''19.52'
$ws.Range("E10").Value = '  +0.74%  '

# Row 11: TRON
$ws.Range("D11").Value = '''0.0842'
$ws.Range("E11").Value = '  +0.45%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '1.824.05'
$ws.Range("E12").Value = '  +0.33%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '1.585.02'
$ws.Range("E13").Value = '  -0.14%  '

# Row 14: Polkadot
$ws.Range("E14").Value = '  +0.47%  '

# Row 15: Polygon
$ws.Range("E15").Value = '  +0.51%  '

# Row 16: Litecoin
$ws.Range("D16").Value = '''65.37'
$ws.Range("E16").Value = '  +1.36%  '

# Row 17: WrappedBTC
$ws.Range("D17").Value = '26.685.42'
$ws.Range("E17").Value = '  +0.32%  '

# Row 18: ShibaInu
$ws.Range("D18").Value = '0.0₃0762'
$ws.Range("E18").Value = '  +4.46%  '

# Row 19: BitcoinCash
$ws.Range("D19").Value = '''210.08'
$ws.Range("E19").Value = '  +1.29%  '

# Row 20: Dai
$ws.Range("E20").Value = '  +0.16%  '

# Row 21: Chainlink
$ws.Range("D21").Value = '''7.16'
$ws.Range("E21").Value = '  +4.51%  '

# Row 22: Uniswap
$ws.Range("E22").Value = '  +0.83%  '

# Row 23: Toncoin
$ws.Range("E23").Value = '  +0.14%  '

# Row 24: Avalanche
$ws.Range("D24").Value = '''8.94'
$ws.Range("E24").Value = '  +0.85%  '

# Row 25: Monero
$ws.Range("D25").Value = '''143.07'
$ws.Range("E25").Value = '  -1.87%  '

# Row 26: BinanceUSD
$ws.Range("E26").Value = '  +0.20%  '

# Row 27: Cosmos
$ws.Range("D27").Value = '''7.14'
$ws.Range("E27").Value = '  +0.10%  '

# Row 28: Stellar
$ws.Range("E28").Value = '  +0.14%  '

# Row 29: EthereumClassic
$ws.Range("D29").Value = '''15.31'
$ws.Range("E29").Value = '  +0.24%  '

# Row 30: Hedera
$ws.Range("D30").Value = '''0.0520'
$ws.Range("E30").Value = '  +3.29%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '  +0.03%  '

# Row 32: Filecoin
$ws.Range("E32").Value = '  +0.35%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = '  +1.87%  '

# Row 34: Maker
$ws.Range("D34").Value = '1.290.34'
$ws.Range("E34").Value = '  +0.70%  '

# Row 35: ImmutableX
$ws.Range("D35").Value = '''0.619'
$ws.Range("E35").Value = '  -5.13%  '

# Row 36: HuobiToken
$ws.Range("E36").Value = '  +0.91%  '

# Row 37: LidoDAOToken
$ws.Range("E37").Value = '  +0.25%  '

# Row 38: VeChain
$ws.Range("E38").Value = '  -0.20%  '

# Row 39: WEMIXToken
$ws.Range("E39").Value = '  +17.15%  '

# Row 40: ARBITRUM
$ws.Range("D40").Value = '''0.828'
$ws.Range("E40").Value = '  -1.19%  '

# Row 41: FraxShare
$ws.Range("E41").Value = '  +0.32%  '

# Row 42: TrustWalletToken -> MXToken (rows 42/43 swapped)
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '''2.19'
$ws.Range("E42").Value = '  -0.47%  '

# Row 43: MXToken -> TrustWalletToken (rows 42/43 swapped)
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''0.783'
$ws.Range("E43").Value = '  -0.48%  '

# Row 44: Aave
$ws.Range("D44").Value = '''63.07'
$ws.Range("E44").Value = '  -0.67%  '

# Row 45: RocketPoolETH
$ws.Range("D45").Value = '1.736.64'
$ws.Range("E45").Value = '  +0.34%  '

# Row 46: Quant
$ws.Range("D46").Value = '''91.36'
$ws.Range("E46").Value = '  +1.85%  '

# Row 47: RenderToken
$ws.Range("E47").Value = '  -1.03%  '

# Row 48: BabyDogeCoin
$ws.Range("E48").Value = '  -0.01%  '

# Row 49: Algorand
$ws.Range("E49").Value = '  -0.23%  '

# Row 50: Cronos
$ws.Range("E50").Value = '  +0.58%  '

# Row 51: USDD
$ws.Range("E51").Value = '  +0.12%  '

Write-Host "Applied cryptos list update"
